$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "5:30 AM"
$ws.Range("C2").Value = 314.296667
$ws.Range("D2").Value = 344.296667
$ws.Range("E2").Value = 13.32
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.333333

$ws.Range("A3").Value = 35
$ws.Range("B3").Value = "5:42 AM"
$ws.Range("C3").Value = 1031.27
$ws.Range("D3").Value = 1061.27
$ws.Range("E3").Value = 2.206667
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

$ws.Range("A4").Value = 84
$ws.Range("B4").Value = "6:07 AM"
$ws.Range("C4").Value = 2509.07
$ws.Range("D4").Value = 2539.07
$ws.Range("E4").Value = 3.404
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.8

$ws.Range("A5").Value = 112
$ws.Range("B5").Value = "6:21 AM"
$ws.Range("C5").Value = 3346.1575
$ws.Range("D5").Value = 3376.1575
$ws.Range("E5").Value = 1.225
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("A6").Value = 132
$ws.Range("B6").Value = "6:31 AM"
$ws.Range("C6").Value = 3944.838333
$ws.Range("D6").Value = 3974.838333
$ws.Range("E6").Value = 1.058333
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0.166667

$ws.Range("A7").Value = 148
$ws.Range("B7").Value = "6:39 AM"
$ws.Range("C7").Value = 4424.9775
$ws.Range("D7").Value = 4454.9775
$ws.Range("E7").Value = 5.47
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0

$ws.Range("A8").Value = 156
$ws.Range("B8").Value = "6:43 AM"
$ws.Range("C8").Value = 4674.78
$ws.Range("D8").Value = 4704.78
$ws.Range("E8").Value = 5.42
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

$ws.Range("A9").Value = 185
$ws.Range("B9").Value = "6:58 AM"
$ws.Range("C9").Value = 5545.74
$ws.Range("D9").Value = 5575.74
$ws.Range("E9").Value = 7.94
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

$ws.Range("A10").Value = 221
$ws.Range("B10").Value = "7:16 AM"
$ws.Range("C10").Value = 6626.76
$ws.Range("D10").Value = 6656.76
$ws.Range("E10").Value = 3.56
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 1

$ws.Range("A11").Value = 249
$ws.Range("B11").Value = "7:29 AM"
$ws.Range("C11").Value = 7455.21
$ws.Range("D11").Value = 7485.21
$ws.Range("E11").Value = 11.435
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0

$ws.Range("A12").Value = 256
$ws.Range("B12").Value = "7:33 AM"
$ws.Range("C12").Value = 7665.623333
$ws.Range("D12").Value = 7695.623333
$ws.Range("E12").Value = 5.7
$ws.Range("F12").Value = 0.346667
$ws.Range("G12").Value = 0

$ws.Range("A13").Value = 271
$ws.Range("B13").Value = "7:41 AM"
$ws.Range("C13").Value = 8128.61
$ws.Range("D13").Value = 8158.61
$ws.Range("E13").Value = 7.54
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

$ws.Range("A14").Value = 281
$ws.Range("B14").Value = "7:45 AM"
$ws.Range("C14").Value = 8407.135
$ws.Range("D14").Value = 8437.135
$ws.Range("E14").Value = 1.65
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0

$ws.Range("A15").Value = 288
$ws.Range("B15").Value = "7:49 AM"
$ws.Range("C15").Value = 8619.565000000001
$ws.Range("D15").Value = 8649.565000000001
$ws.Range("E15").Value = 1.59
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0.5

$ws.Range("A16").Value = 331
$ws.Range("B16").Value = "8:10 AM"
$ws.Range("C16").Value = 9903.9
$ws.Range("D16").Value = 9933.9
$ws.Range("E16").Value = 66.37
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0

$ws.Range("A17").Value = 342
$ws.Range("B17").Value = "8:16 AM"
$ws.Range("C17").Value = 10238.38
$ws.Range("D17").Value = 10268.38
$ws.Range("E17").Value = 2.085
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0

$ws.Range("A18").Value = 351
$ws.Range("B18").Value = "8:20 AM"
$ws.Range("C18").Value = 10514.064286
$ws.Range("D18").Value = 10544.064286
$ws.Range("E18").Value = 2.311429
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0

$ws.Range("A19").Value = 371
$ws.Range("B19").Value = "8:30 AM"
$ws.Range("C19").Value = 11104.845
$ws.Range("D19").Value = 11134.845
$ws.Range("E19").Value = 4.43
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0.5

$ws.Range("A20").Value = 381
$ws.Range("B20").Value = "8:35 AM"
$ws.Range("C20").Value = 11413.92
$ws.Range("D20").Value = 11443.92
$ws.Range("E20").Value = 1.94
$ws.Range("F20").Value = 19.395
$ws.Range("G20").Value = 0

$ws.Range("A21").Value = 488
$ws.Range("B21").Value = "9:29 AM"
$ws.Range("C21").Value = 14621.995
$ws.Range("D21").Value = 14651.995
$ws.Range("E21").Value = 1.415
$ws.Range("F21").Value = 113.915
$ws.Range("G21").Value = 0

$ws.Range("A22").Value = 519
$ws.Range("B22").Value = "9:44 AM"
$ws.Range("C22").Value = 15550.62
$ws.Range("D22").Value = 15580.62
$ws.Range("E22").Value = 2.52
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0

$ws.Range("A23").Value = 526
$ws.Range("B23").Value = "9:48 AM"
$ws.Range("C23").Value = 15765.516667
$ws.Range("D23").Value = 15795.516667
$ws.Range("E23").Value = 5.23
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0

$ws.Range("A24").Value = 545
$ws.Range("B24").Value = "9:57 AM"
$ws.Range("C24").Value = 16330.63
$ws.Range("D24").Value = 16360.63
$ws.Range("E24").Value = 1.686667
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0

$ws.Range("A25").Value = 572
$ws.Range("B25").Value = "10:11 A"
$ws.Range("C25").Value = 17146.84
$ws.Range("D25").Value = 17176.84
$ws.Range("E25").Value = 6.193333
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 1

$ws.Range("A26").Value = 991
$ws.Range("B26").Value = "1:40 PM"
$ws.Range("C26").Value = 29712.875
$ws.Range("D26").Value = 29742.875
$ws.Range("E26").Value = 10.58
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0

$ws.Range("A27").Value = 1000
$ws.Range("B27").Value = "1:45 PM"
$ws.Range("C27").Value = 29983.822
$ws.Range("D27").Value = 30013.822
$ws.Range("E27").Value = 2.52
$ws.Range("F27").Value = 0.482
$ws.Range("G27").Value = 0

$ws.Range("A28").Value = 1042
$ws.Range("B28").Value = "2:06 PM"
$ws.Range("C28").Value = 31251.5775
$ws.Range("D28").Value = 31281.5775
$ws.Range("E28").Value = 11.14
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0.25

$ws.Range("A29").Value = 1065
$ws.Range("B29").Value = "2:17 PM"
$ws.Range("C29").Value = 31927.45
$ws.Range("D29").Value = 31957.45
$ws.Range("E29").Value = 12.815
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0

$ws.Range("A30").Value = 1076
$ws.Range("B30").Value = "2:23 PM"
$ws.Range("C30").Value = 32264.416667
$ws.Range("D30").Value = 32294.416667
$ws.Range("E30").Value = 0.07000000000000001
$ws.Range("F30").Value = 0.333333
$ws.Range("G30").Value = 0

$ws.Range("A31").Value = 1103
$ws.Range("B31").Value = "2:36 PM"
$ws.Range("C31").Value = 33070.88
$ws.Range("D31").Value = 33100.88
$ws.Range("E31").Value = 25.1
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 1

$ws.Range("A32").Value = 1144
$ws.Range("B32").Value = "2:57 PM"
$ws.Range("C32").Value = 34295.06
$ws.Range("D32").Value = 34325.06
$ws.Range("E32").Value = 2.785
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0

$ws.Range("A33").Value = 1154
$ws.Range("B33").Value = "3:02 PM"
$ws.Range("C33").Value = 34618.26
$ws.Range("D33").Value = 34648.26
$ws.Range("E33").Value = 1.505
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0

$ws.Range("A34").Value = 1173
$ws.Range("B34").Value = "3:11 PM"
$ws.Range("C34").Value = 35173.265
$ws.Range("D34").Value = 35203.265
$ws.Range("E34").Value = 2.2825
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0

$ws.Range("A35").Value = 1216
$ws.Range("B35").Value = "3:33 PM"
$ws.Range("C35").Value = 36464.233333
$ws.Range("D35").Value = 36494.233333
$ws.Range("E35").Value = 7.393333
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 1

$ws.Range("A36").Value = 1224
$ws.Range("B36").Value = "3:37 PM"
$ws.Range("C36").Value = 36710.185
$ws.Range("D36").Value = 36740.185
$ws.Range("E36").Value = 3.49
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0.25

$ws.Range("A37").Value = 1235
$ws.Range("B37").Value = "3:42 PM"
$ws.Range("C37").Value = 37028.365
$ws.Range("D37").Value = 37058.365
$ws.Range("E37").Value = 6.5
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0.5

$ws.Range("A38").Value = 1245
$ws.Range("B38").Value = "3:47 PM"
$ws.Range("C38").Value = 37332.09
$ws.Range("D38").Value = 37362.09
$ws.Range("E38").Value = 193.12
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 6

$ws.Range("A39").Value = 1277
$ws.Range("B39").Value = "4:03 PM"
$ws.Range("C39").Value = 38285.38
$ws.Range("D39").Value = 38315.38
$ws.Range("E39").Value = 25.395
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0.5

$ws.Range("A40").Value = 1289
$ws.Range("B40").Value = "4:09 PM"
$ws.Range("C40").Value = 38650.62
$ws.Range("D40").Value = 38680.62
$ws.Range("E40").Value = 88.59
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0

$ws.Range("A41").Value = 1297
$ws.Range("B41").Value = "4:13 PM"
$ws.Range("C41").Value = 38885.46
$ws.Range("D41").Value = 38915.46
$ws.Range("E41").Value = 73.44
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0

$ws.Range("A42").Value = 1305
$ws.Range("B42").Value = "4:18 PM"
$ws.Range("C42").Value = 39149.33
$ws.Range("D42").Value = 39179.33
$ws.Range("E42").Value = 144.36
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0

$ws.Range("A43").Value = 1330
$ws.Range("B43").Value = "4:30 PM"
$ws.Range("C43").Value = 39870.99
$ws.Range("D43").Value = 39900.99
$ws.Range("E43").Value = 42.65
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 6

$ws.Range("A44").Value = 1339
$ws.Range("B44").Value = "4:34 PM"
$ws.Range("C44").Value = 40148.25
$ws.Range("D44").Value = 40178.25
$ws.Range("E44").Value = 8.445
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0

$ws.Range("A45").Value = 1353
$ws.Range("B45").Value = "4:41 PM"
$ws.Range("C45").Value = 40571.0325
$ws.Range("D45").Value = 40601.0325
$ws.Range("E45").Value = 4.135
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0.25

$ws.Range("A46").Value = 1382
$ws.Range("B46").Value = "4:56 PM"
$ws.Range("C46").Value = 41437.05
$ws.Range("D46").Value = 41467.05
$ws.Range("E46").Value = 13.86
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0.5

$ws.Range("A47").Value = 1418
$ws.Range("B47").Value = "5:14 PM"
$ws.Range("C47").Value = 42521.595
$ws.Range("D47").Value = 42551.595
$ws.Range("E47").Value = 25.965
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0

$ws.Range("A48").Value = 1441
$ws.Range("B48").Value = "5:25 PM"
$ws.Range("C48").Value = 43223.54
$ws.Range("D48").Value = 43253.54
$ws.Range("E48").Value = 2.86
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0

$ws.Range("A49").Value = 1471
$ws.Range("B49").Value = "5:41 PM"
$ws.Range("C49").Value = 44128.76
$ws.Range("D49").Value = 44158.76
$ws.Range("E49").Value = 5.15
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 1

$ws.Range("A50").Value = 1481
$ws.Range("B50").Value = "5:45 PM"
$ws.Range("C50").Value = 44405.94
$ws.Range("D50").Value = 44435.94
$ws.Range("E50").Value = 5.48
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 1

$ws.Range("A51").Value = 1518
$ws.Range("B51").Value = "6:04 PM"
$ws.Range("C51").Value = 45519.58
$ws.Range("D51").Value = 45549.58
$ws.Range("E51").Value = 5.9
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0.333333
